# Update the "bigmac_index.py" code sheet to match the latest xlwings API:
# the single import line
#     from xlwings.reports import create_report  # part of xlwings PRO
# is replaced by two lines:
#     # Requires a license key: https://www.xlwings.org/trial
#     from xlwings.pro.reports import create_report

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bigmac_index.py")

# Make room for the extra line by pushing row 5 (and everything below it)
# down by one row.
$ws.Rows.Item(5).Insert()

# Row 4 keeps the "license key" comment, row 5 gets the new import line.
$ws.Cells.Item(4, 1).Value = "# Requires a license key: https://www.xlwings.org/trial"
$ws.Cells.Item(5, 1).Value = "from xlwings.pro.reports import create_report"
